# Update LR-pairs data (Wnt7b-Fzd8) with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5854969999999999
$ws.Range("H2").Value = 1.170994
$ws.Range("M2").Value = 3.1851815
$ws.Range("N2").Value = 6.370363
$ws.Range("O2").Value = 0.4406530230187619
$ws.Range("P2").Value = 0.3851702893788179
$ws.Range("Q2").Value = 1.8649142127055
$ws.Range("R2").Value = 7.459656850821999
$ws.Range("S2").Value = 0.4406530230187619
$ws.Range("T2").Value = 0.3851702893788179

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5854969999999999
$ws.Range("H3").Value = 1.170994
$ws.Range("O3").Value = 0.2827048402157753
$ws.Range("P3").Value = 0.3706641033643825
$ws.Range("Q3").Value = 1.196452190222667
$ws.Range("R3").Value = 7.178713141335999
$ws.Range("S3").Value = 0.2827048402157753
$ws.Range("T3").Value = 0.3706641033643825

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5854969999999999
$ws.Range("H4").Value = 1.170994
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.003928
$ws.Range("N4").Value = 0.011784
$ws.Range("O4").Value = 0.0005434180358066555
$ws.Range("P4").Value = 0.0007124942001013113
$ws.Range("Q4").Value = 0.002299832216
$ws.Range("R4").Value = 0.013798993296
$ws.Range("S4").Value = 0.0005434180358066555
$ws.Range("T4").Value = 0.0007124942001013113

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5854969999999999
$ws.Range("H5").Value = 1.170994
$ws.Range("M5").Value = 1.9606995
$ws.Range("N5").Value = 3.921399
$ws.Range("O5").Value = 0.2712524111754306
$ws.Range("P5").Value = 0.2370989514411984
$ws.Range("Q5").Value = 1.1479836751515
$ws.Range("R5").Value = 4.591934700605999
$ws.Range("S5").Value = 0.2712524111754306
$ws.Range("T5").Value = 0.2370989514411984

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5854969999999999
$ws.Range("H6").Value = 1.170994
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02851766666666667
$ws.Range("N6").Value = 0.085553
$ws.Range("O6").Value = 0.003945268433245655
$ws.Range("P6").Value = 0.005172778029639129
$ws.Range("Q6").Value = 0.01669700828033333
$ws.Range("R6").Value = 0.100182049682
$ws.Range("S6").Value = 0.003945268433245655
$ws.Range("T6").Value = 0.005172778029639129

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5854969999999999
$ws.Range("H7").Value = 1.170994
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.006513000000000001
$ws.Range("N7").Value = 0.019539
$ws.Range("O7").Value = 0.0009010391209798237
$ws.Range("P7").Value = 0.001181383585860448
$ws.Range("Q7").Value = 0.003813341961
$ws.Range("R7").Value = 0.022880051766
$ws.Range("S7").Value = 0.0009010391209798237
$ws.Range("T7").Value = 0.001181383585860448
